$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows (columns A:G only, which is the full extent of the
# table) above the current row 24 (the "(400, 160, 200)" block), shifting
# the existing rows 24-65 down to 26-67.
$ws.Range("A24:G25").Insert(-4121) | Out-Null

# Copy the formatting of the row immediately below (now row 26, the old row 24)
# onto the freshly inserted rows so style/number formats match the rest of the table.
$ws.Range("A26:G26").Copy()
$ws.Range("A24:G25").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# New data rows for the (300, 160, 200) size, inserted between
# (250, 160, 200) [row 23] and (400, 160, 200) [now row 26].
$ws.Range("A24").Value = "(300, 160, 200)"
$ws.Range("B24").Value = "(200, 200, 300)"
$ws.Range("C24").Value = "(300, 160, 200, 300)"
$ws.Range("D24").Value = "([-1],[0])"
$ws.Range("E24").Value = "abc * cde -> abde"
$ws.Range("F24").Value = "float32"
$ws.Range("G24").Value = "abc * cde -> abde-fp32"

$ws.Range("A25").Value = "(300, 160, 200)"
$ws.Range("B25").Value = "(200, 200, 300)"
$ws.Range("C25").Value = "(300, 160, 200, 300)"
$ws.Range("D25").Value = "([-1],[0])"
$ws.Range("E25").Value = "abc * cde -> abde"
$ws.Range("F25").Value = "float16"
$ws.Range("G25").Value = "abc * cde -> abde-fp16"

$ws.Range("E21").Select()
